$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.61
$wsSummary.Range("B4").Value = -0.39
$wsSummary.Range("B5").Value = -0.34
$wsSummary.Range("B6").Value = 23
$wsSummary.Range("B8").Value = 11
$wsSummary.Range("B9").Value = 26.09

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.61
$wsStatus.Range("D4").Value = 23
$wsStatus.Range("E4").Value = -0.39
$wsStatus.Range("F4").Value = -0.39
$wsStatus.Range("G4").Value = 26.09

# ---- New trade row (#23) appended to "All Trades" and "MarketMaking" sheets ----
function Add-TradeRow($ws) {
    $ws.Range("A24").Value = 23

    # Date/Time columns must stay plain text (match formatting of existing rows),
    # so force text format before assignment, then clear the style back to Normal
    # so no extra cell style is left behind.
    $ws.Range("B24").NumberFormat = "@"
    $ws.Range("B24").Value = "2026-02-17"
    $ws.Range("B24").Style = "Normal"

    $ws.Range("C24").NumberFormat = "@"
    $ws.Range("C24").Value = "08:02:51"
    $ws.Range("C24").Style = "Normal"

    $ws.Range("D24").Value = "MarketMaking"
    $ws.Range("E24").Value = "UP"
    $ws.Range("F24").Value = 0.6
    $ws.Range("G24").Value = 0.55
    $ws.Range("H24").Value = "CLOSED"
    $ws.Range("I24").Value = -8.333299999999999
    $ws.Range("J24").Value = -0.05
    $ws.Range("K24").Value = 99.61
    $ws.Range("L24").Value = 0
    $ws.Range("M24").Value = 0
    $ws.Range("N24").Value = 0.6
    $ws.Range("O24").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P24").Value = "early_exit"
    $ws.Range("Q24").Value = 0.13
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
